$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New date labels for the newly appended rows 131-133 (column A),
# shared strings "09 06 2020", "10 06 2020", "11 06 2020" get added
# automatically when assigned since they are not yet present in the table.
$ws.Range("A131").Value = "09 06 2020"
$ws.Range("A132").Value = "10 06 2020"
$ws.Range("A133").Value = "11 06 2020"

# Updated / newly-added numeric data points across rows 127-131.
$cellValues = @{
    "S127" = 18.2865828
    "AM127" = 14.4951626
    "F128" = 15.7906661
    "G128" = 11.47303
    "H128" = 14.5069889
    "L128" = 11.7736866
    "R128" = 15.130351
    "S128" = 18.6131636
    "AJ128" = 15.475009
    "AM128" = 14.5435281
    "AQ128" = 11.7189437
    "AR128" = 16.6666667
    "AT128" = 14.5069622
    "AX128" = 17.8924122
    "B129" = 13.3966891
    "C129" = 19.551786
    "D129" = 17.8413619
    "F129" = 16.2745311
    "G129" = 11.3958285
    "H129" = 14.0148556
    "I129" = 14.5096508
    "J129" = 16.1547912
    "K129" = 14.645359
    "L129" = 11.7526909
    "M129" = 15.6578749
    "O129" = 7.8093812
    "P129" = 20.2247116
    "Q129" = 14.2099706
    "R129" = 14.8746526
    "S129" = 19.0197771
    "T129" = 14.5482429
    "U129" = 14.7450354
    "V129" = 19.0637298
    "W129" = 15.6230926
    "X129" = 15.5382182
    "Y129" = 10.2209302
    "Z129" = 13.3492913
    "AA129" = 16.3383757
    "AB129" = 14.8020335
    "AD129" = 20.4561248
    "AE129" = 9.7618005
    "AF129" = 14.3768644
    "AG129" = 18.7697725
    "AH129" = 21.1396537
    "AI129" = 12.1983914
    "AJ129" = 15.1161879
    "AK129" = 13.332979
    "AL129" = 13.660548
    "AM129" = 14.2235266
    "AN129" = 13.6231201
    "AO129" = 14.0414219
    "AP129" = 11.4422502
    "AQ129" = 11.7920549
    "AS129" = 13.7175911
    "AT129" = 14.4870187
    "AU129" = 19.5797247
    "AV129" = 13.6817435
    "AW129" = 14.6527513
    "AX129" = 17.8005621
    "AY129" = 14.5569006
    "BA129" = 8.9254684
    "BB129" = 12.6542692
    "BC129" = 13.4698016
    "BD129" = 13.9001569
    "BE129" = 15.3266506
    "B130" = 13.9116719
    "C130" = 19.3275917
    "D130" = 18.2124461
    "F130" = 16.4830531
    "G130" = 11.348569
    "H130" = 14.0403854
    "I130" = 14.4266132
    "J130" = 16.581306
    "K130" = 15.0073024
    "L130" = 11.5862401
    "M130" = 15.6111176
    "O130" = 8.5285132
    "P130" = 20.2588454
    "Q130" = 13.8997025
    "R130" = 14.6645239
    "S130" = 19.0685443
    "T130" = 14.8131428
    "U130" = 15.2412698
    "V130" = 18.4916255
    "W130" = 15.3152684
    "X130" = 15.2524546
    "Y130" = 10.3298838
    "Z130" = 13.1744174
    "AA130" = 16.1353727
    "AB130" = 14.7477796
    "AD130" = 21.1474806
    "AE130" = 9.9432611
    "AF130" = 14.6895306
    "AG130" = 17.9764536
    "AH130" = 20.9150219
    "AI130" = 12.297391
    "AJ130" = 15.0489228
    "AK130" = 13.6308568
    "AL130" = 13.6533383
    "AM130" = 13.9227634
    "AN130" = 13.7584017
    "AO130" = 14.0260367
    "AP130" = 11.3501403
    "AQ130" = 11.8950669
    "AS130" = 13.3378153
    "AT130" = 14.5205048
    "AU130" = 20.0684195
    "AV130" = 13.6226858
    "AW130" = 14.9484815
    "AX130" = 17.9805388
    "AY130" = 14.5202266
    "BA130" = 8.3381668
    "BB130" = 12.3891171
    "BC130" = 13.3228056
    "BD130" = 13.0853577
    "BE130" = 14.9665782
    "B131" = 14.4913628
    "C131" = 19.9402313
    "D131" = 18.1086032
    "F131" = 16.5476203
    "G131" = 11.2962349
    "H131" = 13.9327142
    "I131" = 14.2798881
    "J131" = 15.625
    "K131" = 14.2111012
    "L131" = 11.7884436
    "M131" = 15.4547589
    "O131" = 8.5633947
    "P131" = 20.313587
    "Q131" = 14.3406705
    "R131" = 14.4654193
    "S131" = 18.8825783
    "T131" = 14.5313118
    "U131" = 15.0245709
    "V131" = 18.2493735
    "W131" = 15.0174044
    "X131" = 15.3453029
    "Y131" = 10.1209982
    "Z131" = 13.1206828
    "AA131" = 15.8758346
    "AB131" = 14.3983524
    "AD131" = 21.2176896
    "AE131" = 10.4555128
    "AF131" = 14.7950494
    "AG131" = 18.3252985
    "AH131" = 21.1554369
    "AI131" = 12.0874698
    "AJ131" = 14.6992655
    "AK131" = 14.220647
    "AL131" = 13.2396025
    "AM131" = 13.3784943
    "AN131" = 13.598012
    "AO131" = 13.9496459
    "AP131" = 11.6229082
    "AQ131" = 11.6127365
    "AS131" = 13.5426632
    "AT131" = 15.0218712
    "AU131" = 19.7156971
    "AV131" = 13.9147307
    "AW131" = 15.1095526
    "AX131" = 17.9030771
    "AY131" = 14.2341501
    "BA131" = 8.6670746
    "BB131" = 12.4524758
    "BC131" = 13.5970998
    "BD131" = 13.9478458
    "BE131" = 15.1991216
}

foreach ($addr in $cellValues.Keys) {
    $ws.Range($addr).Value = $cellValues[$addr]
}
